$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 8 (3k / 0603 group): R18, R19 -> R1, R18, R19 (new qty 3) ---
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "R1, R18, R19"

# --- Add new BoM row 15 for the indication LED ---
# Insert a fresh row at 15 so it inherits formatting (styles) from the row above it (row 14),
# matching the number formats used throughout the table (General / Text).
$ws.Rows("15:15").Insert(-4121)

$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "LED1"
$ws.Range("C15").Value = "LED_0603"
$ws.Range("E15").Value = "C72043"
$ws.Range("D15").Value = "GREEN"

# --- Grow the query table / list object to include the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E15"))

# --- Keep the ExternalData_1 defined name (query spill range) in sync ---
$extName = $wb.Names.Item("ExternalData_1")
$extName.RefersTo = "=Sheet1!`$A`$1:`$D`$15"

# --- Update the active selection to mirror the post-edit workbook state ---
$ws.Range("D16").Select()
